# Update vm_pu.xlsx res_bus data for Case_1_22 (380 kV case)
# Columns in order correspond to: B,C,D,E,F,I,J,K,L,M,N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(2,1.02,1.042208162133836,1.056128452003317,1.05100371770582,1.06345875501102,1.04168007526778,1.047285222020765,1.058867013043492,1.053756424630235,1.066177338631001,1.019502640262109),
    @(3,1.02,1.043617184785783,1.056901014654519,1.0521782419878,1.064532641303373,1.041865039983536,1.048338304593187,1.059453258015499,1.054742576841256,1.067065590087416,1.019864806263668),
    @(4,1.02,1.044528572998245,1.057400078921198,1.052938043373812,1.065226981497593,1.041983084226765,1.049018952061019,1.059831103004246,1.055379914884067,1.067639213772511,1.020098571780088),
    @(5,1.02,1.044911642639965,1.057609685705385,1.053257419756746,1.065518755291233,1.042032317432823,1.049304915301941,1.059989591587463,1.055647670061635,1.06788009501459,1.020196708735529),
    @(6,1.02,1.044975957271748,1.057644867893148,1.05331104193948,1.065567738001627,1.042040560885592,1.049352919257092,1.060016181567974,1.055692616741319,1.067920524234222,1.020213178289302),
    @(7,1.02,1.044533691894526,1.057402880482343,1.052942311067046,1.065230880689789,1.041983743625793,1.049022773823452,1.05983322214054,1.055383493353212,1.067642433499241,1.020099883632396),
    @(8,1.02,1.04268442169057,1.056389716255402,1.051400695384018,1.063821791813216,1.04174292465306,1.047641276758001,1.059065447001678,1.054089859821612,1.066477763114772,1.01962515656132),
    @(9,1.02,1.039422946078511,1.054597984343518,1.048682574495007,1.061334623426622,1.041306003379611,1.045200901440596,1.057701063238158,1.051804327653391,1.064416712436578,1.018784147363504),
    @(10,1.02,1.037246457821375,1.053399172848079,1.046869265996361,1.059673606361047,1.041006264684498,1.043569774849578,1.0567837367134,1.050276474224691,1.063036700887827,1.018220410763449),
    @(11,1.02,1.036303438479407,1.052879044115345,1.046083757115831,1.058953655733594,1.040874466122185,1.042862441011903,1.056384681783918,1.049613876678494,1.062437701238899,1.017975568001989),
    @(12,1.02,1.035953065811269,1.052685689057571,1.045791930802335,1.0586861238621,1.040825208060396,1.042599545203765,1.056236176866537,1.049367601300781,1.062214986953981,1.017884510090795),
    @(13,1.02,1.03602822624551,1.052727171445687,1.045854531001103,1.058743515375056,1.040835787761743,1.042655944567394,1.056268044274158,1.049420435362655,1.062262769872649,1.017904047425749),
    @(14,1.02,1.036274478499442,1.052863064513178,1.046059635767903,1.058931543724864,1.040870400601941,1.042840713257325,1.056372411999561,1.049593522701307,1.06241929608856,1.017968043425927),
    @(15,1.02,1.036426190118156,1.052946772018035,1.046186000422569,1.059047379537819,1.040891686661627,1.042954533969338,1.056436679520178,1.049700146613874,1.062515707958235,1.018007458546594),
    @(16,1.02,1.037309030118397,1.053433670219065,1.046921390324598,1.059721371780784,1.041014969355538,1.043616695947045,1.056810181708981,1.050320426761433,1.063076423898184,1.018236644466239),
    @(17,1.02,1.037862651802995,1.053738811013344,1.047382589279412,1.060143955266796,1.041091763082717,1.044031770357872,1.057043974891296,1.050709235155899,1.063427757795333,1.018380207575392),
    @(18,1.02,1.038185513828137,1.053916694608106,1.047651566780484,1.060390371808048,1.041136361673677,1.044273775624118,1.057180164390717,1.050935921625561,1.063632545574457,1.018463874076085),
    @(19,1.02,1.038295592073703,1.053977331392818,1.047743275815134,1.060474381796001,1.04115153576126,1.044356276204782,1.057226571301948,1.05101319916635,1.06370234933954,1.018492390106609),
    @(20,1.02,1.03780325928224,1.053706082628289,1.047333110334718,1.060098623237767,1.041083543899857,1.043987247235327,1.057018909523995,1.050667529945136,1.063390077420864,1.018364812009132),
    @(21,1.02,1.036201965966631,1.052823051704912,1.045999239023707,1.058876177119943,1.040860216327789,1.042786307942801,1.056341685989949,1.049542557155251,1.062373209084642,1.017949201310545),
    @(22,1.02,1.035194626867442,1.052266951770511,1.04516027102506,1.058106938371062,1.040718052803937,1.042030299246622,1.055914279020284,1.048834331480079,1.061732594712069,1.0176872394779),
    @(23,1.02,1.035728689581891,1.052561836672432,1.045605054277166,1.058514787596272,1.040793582216971,1.042431163129365,1.056141008366307,1.049209862587039,1.062072317448479,1.017826172456503),
    @(24,1.02,1.037830096379005,1.053720871485429,1.047355467838961,1.060119107044398,1.041087258393935,1.044007365628659,1.057030236034953,1.050686375040154,1.063407103990681,1.018371768823515),
    @(25,1.02,1.040266478417099,1.055061950835986,1.049385478977349,1.061978120782451,1.041420448136149,1.045832526071394,1.058055150784177,1.052395916241897,1.06495059019789,1.019002104578455)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $val = $rowData[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}
